$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the literal text into the cell (leading apostrophe = text marker so
    # Excel doesn't silently reinterpret plain-decimal-looking strings, e.g.
    # "239.60", as a Number and drop the trailing zero). ClearFormats() then
    # strips the quotePrefix/number-format style Excel stamps on the cell so
    # the cell keeps its original (default) style index.
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "91.292.40"
Set-TextValue "E2" "  +1.90%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.155.43"
Set-TextValue "E3" "  +2.70%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.05%  "

# Row 5 - Solana
Set-TextValue "D5" "239.60"
Set-TextValue "E5" "  +2.02%  "

# Row 6 - BNB
Set-TextValue "D6" "619.19"
Set-TextValue "E6" "  +0.27%  "

# Row 7 - XRP
Set-TextValue "E7" "  +5.28%  "

# Row 8 - Dogecoin
Set-TextValue "E8" "  +3.78%  "

# Row 9 - USDC
Set-TextValue "E9" "  +0.01%  "

# Row 10 - LidoStakedEther
Set-TextValue "D10" "3.153.81"
Set-TextValue "E10" "  +2.70%  "

# Row 11 - Cardano
Set-TextValue "D11" "0.744"
Set-TextValue "E11" "  +5.32%  "

# Row 12 - TRON
Set-TextValue "E12" "  +2.17%  "

# Row 13 - ShibaInu
Set-TextValue "E13" "  +0.90%  "

# Row 14 - Avalanche
Set-TextValue "D14" "35.19"
Set-TextValue "E14" "  +1.17%  "

# Row 15 - Toncoin
Set-TextValue "E15" "  +4.47%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "91.075.20"
Set-TextValue "E16" "  +1.92%  "

# Row 17 - WrappedliquidstakedEther2.0
Set-TextValue "D17" "3.746.15"

# Row 18 - WrappedEther
Set-TextValue "D18" "3.189.87"
Set-TextValue "E18" "  +3.82%  "

# Row 19 - SuiNetwork
Set-TextValue "D19" "3.75"
Set-TextValue "E19" "  +0.32%  "

# Row 20 - Chainlink
Set-TextValue "E20" "  +11.31%  "

# Row 21 - Polkadot
Set-TextValue "D21" "6.05"
Set-TextValue "E21" "  +12.60%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "454.10"
Set-TextValue "E22" "  +5.72%  "

# Row 23 - PEPE
Set-TextValue "E23" "  -3.85%  "

# Row 24 - Uniswap
Set-TextValue "E24" "  +6.35%  "

# Row 25 - NEARProtocol
Set-TextValue "D25" "6.00"
Set-TextValue "E25" "  +8.33%  "

# Row 26 - Litecoin
Set-TextValue "D26" "88.94"
Set-TextValue "E26" "  +2.13%  "

# Row 27 - Aptos
Set-TextValue "E27" "  +3.56%  "

# Row 29 - Dai
Set-TextValue "D29" "0.999"
Set-TextValue "E29" "  -0.19%  "

# Row 30 - Hedera
Set-TextValue "D30" "0.131"
Set-TextValue "E30" "  +45.64%  "

# Row 31 - Stellar
Set-TextValue "D31" "0.237"
Set-TextValue "E31" "  +18.85%  "

# Row 32 - Cronos
Set-TextValue "D32" "0.172"
Set-TextValue "E32" "  +10.93%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue "E33" "  +4.85%  "

# Row 34 - Kaspa
Set-TextValue "E34" "  +15.11%  "

# Row 35 - Binance-PegBSC-USD
Set-TextValue "E35" "  -5.31%  "

# Row 36 - RenderToken
Set-TextValue "D36" "7.65"
Set-TextValue "E36" "  +9.70%  "

# Row 37 - EthereumClassic
Set-TextValue "D37" "26.51"
Set-TextValue "E37" "  +3.71%  "

# Row 38 - Bittensor
Set-TextValue "D38" "511.60"
Set-TextValue "E38" "  +4.50%  "

# Row 39 - PancakeSwap
Set-TextValue "E39" "  +5.03%  "

# Row 40 - Fetch.AI
Set-TextValue "D40" "1.35"
Set-TextValue "E40" "  +8.72%  "

# Row 41 - MantraDAO
Set-TextValue "D41" "3.88"
Set-TextValue "E41" "  -3.39%  "

# Row 42 - PolygonEcosystemToken
Set-TextValue "D42" "0.450"
Set-TextValue "E42" "  +13.56%  "

# Row 43 - dogwifhat
Set-TextValue "D43" "3.48"
Set-TextValue "E43" "  -3.07%  "

# Row 44 - WhiteBITCoin
Set-TextValue "E44" "  +0.34%  "

# Row 45 - USDe
Set-TextValue "E45" "  -0.02%  "

# Row 46 - ARBITRUM
Set-TextValue "D46" "0.723"
Set-TextValue "E46" "  +7.45%  "

# Row 47 - Stacks
Set-TextValue "E47" "  +5.94%  "

# Row 48 - Monero
Set-TextValue "D48" "156.62"
Set-TextValue "E48" "  -0.44%  "

# Row 49 - ImmutableX
Set-TextValue "E49" "  +6.97%  "

# Row 50 - Filecoin
Set-TextValue "E50" "  +4.88%  "

# Row 51 - OKB -> VeChain
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D51" "0.0327"
Set-TextValue "E51" "  +14.87%  "
